$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new edge rows (OBJECTID, START WATER NODE ID, END WATER NODE ID)
$newRows = @(
    @(71, 3, 11),
    @(72, 11, 13),
    @(73, 9, 14),
    @(74, 13, 43)
)

$startRow = 72
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]
    $rowRange = $ws.Range("A" + $r + ":C" + $r)
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $rowRange.HorizontalAlignment = -4108
    $rowRange.VerticalAlignment = -4108
}

# Update selection to match the new last row
$ws.Range("A75").Select()
